# "transport limit h2 pipeline implemented"
#
# Adds a new "MWh/y" column (F) to the "2045_level_2" sheet that converts
# the existing GWh/d transport-limit column (E) into yearly MWh
# (E * 1000 * 365), and labels the corresponding unit ("GWh/d") on the
# "tidy" sheet's helper header row so the lookup table stays self
# documenting.

$wb = $excel.ActiveWorkbook

# --- "2045_level_2" sheet: new column F = MWh/y -----------------------
$ws1 = $wb.Worksheets.Item("2045_level_2")
$ws1.Activate()

$ws1.Range("F1").Value = "MWh/y"

for ($r = 2; $r -le 85; $r++) {
    $ws1.Cells.Item($r, 6).Formula = "=E$r*1000*365"
}

# Matches the author's final selection after filling the new column down.
$ws1.Range("F2:F85").Select()

# --- "tidy" sheet: note the unit used by the lookup table -------------
$ws2 = $wb.Worksheets.Item("tidy")
$ws2.Activate()

$ws2.Range("P1").Value = "GWh/d"

$ws2.Range("P3").Select()

# Leave "2045_level_2" as the active/visible tab, matching tabSelected.
$ws1.Activate()
